$d = $word.ActiveDocument

function Replace-Text($oldText, $newText, $label) {
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    Write-Output ($label + ": " + $found)
    if (-not $found) {
        throw ("Could not find text for: " + $label)
    }
}

# Heading3 title change: Mechanical Properties -> Physical Metallurgy
Replace-Text "Mechanical Properties" "Physical Metallurgy" "Replace 0 (heading3)"

# Ativacao date change: 2024 -> 2025
Replace-Text "Ativação: 01/01/2024" "Ativação: 01/01/2025" "Replace 1 (ativacao)"

# Objetivos paragraph full text change
$old2 = "Esta disciplina faz parte da formação do engenheiro de materiais e têm como objetivo gerar competências no desenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de materiais e a redução de ocorrência de falhas estruturais. Para tanto, a disciplina estabelece correlações com outras do curso de Engenharia de Materiais como LOM3013 – Ciência dos Materiais, LOM3057 – Introdução aos Materiais Poliméricos, LOM3032 - Cerâmica Física e LOM3011- Ensaios Mecânicos. Desta forma, são apresentadas a correlação entre propriedades e microestrutura de materiais para aplicações em Engenharia permitindo aos alunos a prática da redação científica e da busca bibliográfica para incentivar a solução de problemas em engenharia."
$new2 = "Esta disciplina faz parte da formação do engenheiro de materiais e têm como objetivo gerar competências nodesenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de metais e a redução de ocorrência de falhas estruturais baseado no trinômio propriedades, estrutura metalúrgica e processamento metalúrgico dos metais aplicado a engenharia permitindo aos alunos a prática da redação científica e da busca de projetos para incentivar a solução de problemas em engenharia."
Replace-Text $old2 $new2 "Replace 2 (objetivos)"

# Insert new docente entry before the "7459752" run (Cassius Olivio ...)
$oldDocente = "7459752 - Maria Ismenia Sodero Toledo Faria"
$rngDocente = $d.Content
$foundDocente = $rngDocente.Find.Execute($oldDocente, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("Found docente target: " + $foundDocente)
if (-not $foundDocente) {
    throw "Could not find docente anchor text"
}
$insertRange = $d.Range($rngDocente.Start, $rngDocente.Start)
$insertRange.InsertBefore("3586455 - Cassius Olivio Figueiredo Terra Ruchert`v")
Write-Output "Inserted docente"

# Programa resumido paragraph full text change
$old3 = "1. Introdução ao conceito de propriedades mecânicas. 2. Elasticidade e Mecanismos de deformação plástica. 3. Teoria das discordâncias. 4.Mecanismos de endurecimento. 5. Comportamento mecânico dos materiais metálicos. 6. Estudo comparativo de propriedades mecânicas de materiais metálicos, cerâmicos e poliméricos. 7. Influência da temperatura no comportamento mecânico de materiais. 8. Introdução básica à análise de falhas de materiais dúcteis e frágeis."
$new3 = "1. Introdução ao conceito de propriedades mecânicas. 2. Elasticidade e mecanismos de deformação plástica. 3. Teoria das discordâncias. 4.Mecanismos de endurecimento. 5. Comportamento mecânico dos materiais metálicos. 6. Estudo comparativo de propriedades mecânicas de materiais metálicos. 7. Influência da temperatura no comportamento mecânico dos metais. 8. Introdução básica à análise de falhas de metais dúcteis e frágeis."
Replace-Text $old3 $new3 "Replace 3 (programa resumido)"

# Programa paragraph full text change
$old4 = "1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas de materiais. Comportamento elástico e plástico de metais e ligas. 2. MECANISMOS DE DEFORMAÇÃO PLÁSTICA: Sistemas de deslizamento e movimentação de discordâncias. Deformação por maclação Movimento relativo de grãos. Difusão. 3. TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Aços comuns e especiais. Tratamentos térmicos em aços. 5. COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência, fadiga de alto ciclo e propagação de trincas por fadiga. Impacto e a transição dúctil-frágil. 6. COMPORTAMENTO MECÂNICO DE MATERIAIS CERÂMICOS E POLIMÉRICOS: Estudo comparativo de propriedades mecânicas de materiais metálicos, cerâmicos e poliméricos 7. Influência da temperatura sobre o comportamento mecânico de materiais. Aspectos básicos  da  análise de falhas em materiais metálicos, cerâmicos e poliméricos."
$new4 = "Programa1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas de materiais metálicos. Comportamento elástico e plástico de metais suas ligas e materiais não ferrosos. 2. MECANISMOS DE DEFORMAÇÃO PLÁSTICA: Sistemas de deslizamento e movimentação de discordâncias. Deformação por maclação. Movimento relativo de grãos. Difusão. 3. TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Aços comuns e especiais. Estudo de ligas não metálicas. Tratamentos térmicos em aços e ligas especiais. 5. COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência. Impacto e a transição dúctil-frágil. 6. Influência da temperatura sobre o comportamento mecânico dos metais. Aspectos básicos da análise de falhas em materiais metálicos."
Replace-Text $old4 $new4 "Replace 4 (programa)"

# Norma de recuperacao text full replacement
$old5 = "1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 2. A. S. Lisbão, Estrutura e propriedades dos polímeros, EduFSCar, São Carlos, 2009. 3. T. H. Courtney, Mechanical Behavior of Materials, Waveland Press, 2005. 4. A. K. Bhargava, Engineering Materials: Polymers, Ceramics and Composites, PHI Learning Pvt. Ltd., 2012. 5.Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 2007. 6. Hull, D. Introduction to Dislocations, Pergamon Press, 1965. 7. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967. 8. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. 9. Van Vlack, L.H. Princípios de Ciência dos Materiais, Ed. Edgard Blucher Ltda., 1970. 10. Costa e Silva, A. L., Mei, P. R. Aços e Ligas especiais, Ed. Edgar Blücher, 2008. 11. Dieter, G.E. Metalurgia Mecânica, Ed. Guanabara Dois, 1986.  12. Callister, W. Ciência e engenharia dos materiais: Uma introdução, Rio de Janeiro, Livros Técnicos e Científicos, 2008. 13. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993."
$new5 = "A recuperação será composta por uma única prova (PR) abrangendo toda a matéria ministrada ao longo do semestre. A Média final (MF) será computada pela relação:  MF = (NF + PR)/2."
Replace-Text $old5 $new5 "Replace 5 (norma de recuperacao)"
